$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Najmuddin" record (row 5) — duplicate/invalid entry removed as
# part of adding unique validation. Remaining rows shift up automatically.
$ws.Rows(5).Delete()

# Append a new student record at the end of the table (now row 13).
$ws.Range("A13").Value = "gksjad"
$ws.Range("B13").Value = 87436252
$ws.Range("C13").Value = "fsdck@mail.sd"
$ws.Range("D13").Value = "Wrestling"
